$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) for 4 rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 368
$ws1.Range("F3").Value = 787
$ws1.Range("F5").Value = 858
$ws1.Range("F6").Value = 2110

# Sheet "全部类型" - update 想去人数 (F column) for the matching 4 rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 368
$ws4.Range("F3").Value = 787
$ws4.Range("F7").Value = 858
$ws4.Range("F8").Value = 2110
